$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Bozoh pass" tipseur row entirely (was row 7), shifting rows 8-17 up by one.
$ws.Rows(7).Delete()

# Re-number the id_requete sequence (column D) so it stays contiguous 1..16.
$ws.Range("D7").Value = 6
$ws.Range("D8").Value = 7
$ws.Range("D9").Value = 8
$ws.Range("D10").Value = 9
$ws.Range("D11").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("D13").Value = 12
$ws.Range("D14").Value = 13
$ws.Range("D15").Value = 14
$ws.Range("D16").Value = 15
$ws.Range("D17").Value = 16

# Update the search query / tweet id info for the remaining rows.
$ws.Range("E2").Value = "Buteurs from:ValueBet3 -is:retweet -is:reply "
$ws.Range("F2").Value = 1561437435409613000

$ws.Range("E4").Value = "has:images has:hashtags -✅ from:Tyldumia -is:reply -is:retweet"
$ws.Range("F4").Value = 1561428521993503000

$ws.Range("E5").Value = ""

$ws.Range("F6").Value = 1561428580210512000

$ws.Range("E7").Value = """si tu suis"" from:BozohPronos -is:reply -is:retweet"
$ws.Range("F7").Value = 1561301884983542000
